# Re-upload of the "Orthopedics United" roster: the player list (column A),
# along with each player's position (column B) and team (column C), is
# reordered. Rows 2-19 are rewritten in place so every player keeps his own
# correct position/team while moving to his new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$roster = @(
    @("Tyus Jones",        "PG",      "Phoenix Suns"),
    @("Jalen Green",       "PG,SG",   "Houston Rockets"),
    @("Keon Johnson",      "PG,SG",   "Brooklyn Nets"),
    @("Chris Paul",        "PG",      "San Antonio Spurs"),
    @("Pascal Siakam",     "SF,PF,C", "Indiana Pacers"),
    @("Draymond Green",    "PF,C",    "Golden State Warriors"),
    @("Khris Middleton",   "SF",      "Milwaukee Bucks"),
    @("Deni Avdija",       "SF,PF",   "Portland Trail Blazers"),
    @("Nikola Jokic",      "C",       "Denver Nuggets"),
    @("Rudy Gobert",       "C",       "Minnesota Timberwolves"),
    @("Jakob Poeltl",      "C",       "Toronto Raptors"),
    @("Jonas Valanciunas", "C",       "Washington Wizards"),
    @("Dejounte Murray",   "PG,SG",   "New Orleans Pelicans"),
    @("Russell Westbrook", "PG,SG",   "Denver Nuggets"),
    @("Jaylen Brown",      "SG,SF",   "Boston Celtics"),
    @("Paolo Banchero",    "SF,PF",   "Orlando Magic"),
    @("Chet Holmgren",     "PF,C",    "Oklahoma City Thunder"),
    @("Jalen Suggs",       "PG,SG",   "Orlando Magic")
)

$row = 2
foreach ($player in $roster) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $row++
}
